$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.353.83'
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Value = '2.011.69'
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.60'
$ws.Range("E5").Value = '  +4.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.05'
$ws.Range("E8").Value = '  -6.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.387'
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0770'
$ws.Range("E10").Value = '  -5.11%  '

$ws.Range("E11").Value = '  -1.75%  '

$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '2.310.60'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.22'
$ws.Range("E13").Value = '  -5.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.801'
$ws.Range("E14").Value = '  -5.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.90'
$ws.Range("E15").Value = '  -6.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.24'
$ws.Range("E16").Value = '  -3.80%  '

$ws.Range("D17").Value = '2.022.14'
$ws.Range("E17").Value = '  +0.95%  '

$ws.Range("D18").Value = '37.268.22'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.79'
$ws.Range("E19").Value = '  -0.65%  '

$ws.Range("D20").Value = '0.0₃0834'
$ws.Range("E20").Value = '  -3.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.15'
$ws.Range("E21").Value = '  -0.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.57'
$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.64'
$ws.Range("E23").Value = '  +6.17%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.72'
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  -5.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.65'
$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.130'
$ws.Range("E29").Value = '  -5.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").Value = '  -5.60%  '

$ws.Range("E31").Value = '  -1.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.64'
$ws.Range("E32").Value = '  -3.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0645'
$ws.Range("E33").Value = '  -1.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  +1.78%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.38'
$ws.Range("E35").Value = '  -2.67%  '

$ws.Range("E36").Value = '  +0.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.27'
$ws.Range("E39").Value = '  -1.77%  '

$ws.Range("E40").Value = '  +4.03%  '

$ws.Range("E41").Value = '  +2.00%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0935'
$ws.Range("E42").Value = '  -4.77%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0212'
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").Value = '1.408.45'
$ws.Range("E44").Value = '  +2.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.94'
$ws.Range("E45").Value = '  -1.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.67'
$ws.Range("E46").Value = '  -5.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.02'
$ws.Range("E47").Value = '  -2.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.06'
$ws.Range("E48").Value = '  -4.36%  '

$ws.Range("E49").Value = '  +2.46%  '

$ws.Range("D50").Value = '2.202.89'
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -7.18%  '
